# Update odds values on the active sheet to reflect the latest FlashScore scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.33
$ws.Range("I2").Value = 1.95
$ws.Range("J2").Value = 5
$ws.Range("L2").Value = 2.75
$ws.Range("Q2").Value = 1.79
$ws.Range("R2").Value = 1.94
$ws.Range("AC2").Value = 9.5
$ws.Range("AD2").Value = 21
$ws.Range("AF2").Value = 51
$ws.Range("AH2").Value = 51
$ws.Range("AN2").Value = 5.5
$ws.Range("AO2").Value = 8
$ws.Range("AP2").Value = 9

# Row 3
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7

# Row 4
$ws.Range("H4").Value = 3.6
$ws.Range("J4").Value = 2.2
$ws.Range("K4").Value = 2.05
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.75
$ws.Range("S4").Value = 2.38
$ws.Range("T4").Value = 1.57
$ws.Range("W4").Value = 4.5
$ws.Range("X4").Value = 1.2
$ws.Range("AA4").Value = 2.5
$ws.Range("AB4").Value = 1.5
$ws.Range("AI4").Value = 7
$ws.Range("AL4").Value = 101
$ws.Range("AN4").Value = 13

# Row 5
$ws.Range("W5").Value = 4.33
$ws.Range("X5").Value = 1.22

# Row 6
$ws.Range("N6").Value = 10
$ws.Range("S6").Value = 2.03
$ws.Range("T6").Value = 1.83
$ws.Range("W6").Value = 3.4
$ws.Range("X6").Value = 1.3
